$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above row 2; the old row 2 becomes row 3, ...,
# the old row 157 becomes row 158 (dimension grows to A1:F158).
$ws.Rows.Item(2).Insert()

# Copy the formatting used by the data rows onto the new (blank) row 2.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)

# The newly inserted row 2 duplicates the row that is now row 3 (the
# former row 2), except it carries the new "Date" value.
$ws.Range("A2").Value = "15-11-2025"
$ws.Range("B2").Value = $ws.Range("B3").Value()
$ws.Range("C2").Value = $ws.Range("C3").Value()
$ws.Range("D2").Value = $ws.Range("D3").Value()
$ws.Range("E2").Value = $ws.Range("E3").Value()
$ws.Range("F2").Value = $ws.Range("F3").Value()

# The row-insert operation moves cell text but does not relocate the
# Hyperlinks collection entries, so what was row 96's link (now row 97)
# needs its own hyperlink object, matching the source workbook's change.
$ws.Hyperlinks.Add($ws.Range("F97"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Adding the hyperlink re-styles the cell with the built-in "Hyperlink"
# look; restore the plain data-row formatting used throughout the sheet.
$ws.Range("F96").Copy()
$ws.Range("F97").PasteSpecial(-4122)
